$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F8").Value = 474
$wsExhibition.Range("F9").Value = 4420
$wsExhibition.Range("F17").Value = 2962
$wsExhibition.Range("F18").Value = 1795
$wsExhibition.Range("F27").Value = 2297
$wsExhibition.Range("F28").Value = 994
$wsExhibition.Range("F29").Value = 2371
$wsExhibition.Range("F32").Value = 550
$wsExhibition.Range("F36").Value = 1091
$wsExhibition.Range("F40").Value = 762

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F9").Value = 474
$wsAll.Range("F10").Value = 4420
$wsAll.Range("F17").Value = 2962
$wsAll.Range("F19").Value = 1795
$wsAll.Range("F29").Value = 2297
$wsAll.Range("F32").Value = 994
$wsAll.Range("F33").Value = 2371
$wsAll.Range("F35").Value = 550
$wsAll.Range("F37").Value = 1091
$wsAll.Range("F40").Value = 762
